# Weekly refresh of the cryptocurrency price / 1h-volume table (GitHub Actions bot).
# Column D (Price) holds numeric-looking text (e.g. "13.00", "62.351.36") in the
# source data, so those cells are forced to Text format before the write so Excel
# does not silently reinterpret them as numbers and drop the formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.351.36'
$ws.Range("E2").Value = '  -3.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.370.24'
$ws.Range("E3").Value = '  -3.96%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.65'
$ws.Range("E5").Value = '  -3.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '124.67'
$ws.Range("E6").Value = '  -7.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.370.44'
$ws.Range("E8").Value = '  -3.91%  '
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("E10").Value = '  -4.83%  '
$ws.Range("E11").Value = '  -4.52%  '
$ws.Range("E12").Value = '  -3.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.953.20'
$ws.Range("E13").Value = '  -3.79%  '
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.377.95'
$ws.Range("E15").Value = '  -3.75%  '
$ws.Range("E16").Value = '  -5.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.403.09'
$ws.Range("E17").Value = '  -3.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.40'
$ws.Range("E18").Value = '  -5.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.16'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.63'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.00'
$ws.Range("E21").Value = '  -4.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '370.34'
$ws.Range("E22").Value = '  -6.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  -4.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.508.78'
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.12'
$ws.Range("E27").Value = '  -10.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.90'
$ws.Range("E29").Value = '  -6.62%  '
$ws.Range("E30").Value = '  -7.18%  '
$ws.Range("E31").Value = '  -6.09%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  -5.32%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.402.22'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("E35").Value = '  -6.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.62'
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.18'
$ws.Range("E37").Value = '  -3.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '164.98'
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.58'
$ws.Range("E39").Value = '  -5.58%  '
$ws.Range("E40").Value = '  -5.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0747'
$ws.Range("E41").Value = '  -5.14%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.763'
$ws.Range("E43").Value = '  -5.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.43'
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.22'
$ws.Range("E45").Value = '  -4.97%  '
$ws.Range("E46").Value = '  -7.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.45'
$ws.Range("E47").Value = '  -9.80%  '
$ws.Range("E48").Value = '  -8.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.56'
$ws.Range("E49").Value = '  -3.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.231.48'
$ws.Range("E50").Value = '  -6.22%  '
$ws.Range("E51").Value = '  -6.86%  '
